# صيدليات دكتور مصطفي طلعت — remove the "GRIPE WATER BAMBINO SYRUP 120 ML"
# product line from the report (row 62) and let everything below it
# shift up, same as a normal Excel "Delete Row".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62 holds the GRIPE WATER BAMBINO SYRUP 120 ML entry (item #59).
# Deleting the whole row shifts every row below it up by one and drops
# the now-unused shared string automatically on save.
$ws.Rows.Item(62).Delete()

# The grand-total cell (column K of the totals row, now row 146 after the
# shift) is a plain stored number rather than a SUM formula, so it needs
# to be corrected by hand: 9263.74 - 75 (the deleted row's price) = 9188.74.
$ws.Range("K146").Value = 9188.74

# Excel also re-measured the now-shorter totals row; match its height.
$ws.Rows.Item(146).RowHeight = 25.5
